$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect (the legacy password hash "D382" isn't
# validated by this host, matching real Excel's "accidental edit guard"
# semantics for sheet protection) so the locked data cells can be written.
$ws.Unprotect("D382")

# --- Update the "as of" date in the confidential disclaimer text (A41) ---
$ws.Range("A41").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-23 for illustrative purposes only and are subject to change."

# --- Refresh the return/percentile figures in columns D (percentile) and E (return) for rows 2-38 ---
$ws.Range("D2").Value = 0.02918399291800076
$ws.Range("E2").Value = -0.01000944287063266
$ws.Range("D3").Value = 0.02877125786888804
$ws.Range("E3").Value = 0.0172188943657301
$ws.Range("D4").Value = 0.03064033076463884
$ws.Range("E4").Value = -0.02391033623910332
$ws.Range("D5").Value = 0.06594581829652141
$ws.Range("E5").Value = 0.00856030628087967
$ws.Range("D6").Value = 0.01410019095620557
$ws.Range("E6").Value = 0.01691347816282041
$ws.Range("D7").Value = 0.01421254307537852
$ws.Range("E7").Value = -0.011708553956298
$ws.Range("D8").Value = 0.03175834034297498
$ws.Range("E8").Value = -0.01862309263486717
$ws.Range("D9").Value = 0.03118301269875726
$ws.Range("E9").Value = -0.004894629503738979
$ws.Range("D10").Value = 0.0332674624871113
$ws.Range("E10").Value = 0.02290149299382538
$ws.Range("D11").Value = 0.03015954848863081
$ws.Range("E11").Value = 0.00231950067476383
$ws.Range("D12").Value = 0.01545795571715418
$ws.Range("E12").Value = -0.04267690619857378
$ws.Range("D13").Value = 0.01653356826938731
$ws.Range("E13").Value = -0.03177167474421116
$ws.Range("D14").Value = 0.01512556303250665
$ws.Range("E14").Value = -0.03094517322569801
$ws.Range("D15").Value = 0.008028300862864275
$ws.Range("E15").Value = -0.02851711026615977
$ws.Range("D16").Value = 0.006875101752862647
$ws.Range("E16").Value = -0.04292057227429702
$ws.Range("D17").Value = 0.03192898837303956
$ws.Range("E17").Value = -0.0007900729655619942
$ws.Range("D18").Value = 0.03211871506485049
$ws.Range("E18").Value = 0.004501234209379978
$ws.Range("D19").Value = 0.03185818533944754
$ws.Range("E19").Value = -0.004990518015769996
$ws.Range("D20").Value = 0.03111305760568731
$ws.Range("E20").Value = -0.009913470055188389
$ws.Range("D21").Value = 0.04321507672165843
$ws.Range("E21").Value = 0.007048989742910683
$ws.Range("D22").Value = 0.02877316573506268
$ws.Range("E22").Value = -0.02227919724162319
$ws.Range("D23").Value = 0.03062549180550279
$ws.Range("E23").Value = 0.003599363189589644
$ws.Range("D24").Value = 0.03022250807239377
$ws.Range("E24").Value = -0.0155293226437726
$ws.Range("D25").Value = 0.01507638248222717
$ws.Range("E25").Value = -0.05624296962879638
$ws.Range("D26").Value = 0.01408323214576437
$ws.Range("E26").Value = -0.05810190411680594
$ws.Range("D27").Value = 0.03044212466760732
$ws.Range("E27").Value = -0.0006893910379165424
$ws.Range("D28").Value = 0.03262535952678135
$ws.Range("E28").Value = -0.02053227986277151
$ws.Range("D29").Value = 0.03001582257014164
$ws.Range("E29").Value = 0.006737573626001137
$ws.Range("D30").Value = 0.02931118399630976
$ws.Range("E30").Value = -0.008317060823027411
$ws.Range("D31").Value = 0.03354346712704184
$ws.Range("E31").Value = -0.008759124087591164
$ws.Range("D32").Value = 0.03213270608346448
$ws.Range("E32").Value = 0.01204644412191591
$ws.Range("D33").Value = 0.02980850111249796
$ws.Range("E33").Value = -0.01954258405871312
$ws.Range("D34").Value = 0.03110754599229392
$ws.Range("E34").Value = 0.004470370168456528
$ws.Range("D35").Value = 0.03108592350898139
$ws.Range("E35").Value = 0.0001159285879896998
$ws.Range("D36").Value = 0.0288030556384653
$ws.Range("E36").Value = -0.01000198714976475
$ws.Range("D37").Value = 0.03086651889889836
$ws.Range("E37").Value = 0.01186069351061403
$ws.Range("D38").Value = 1
$ws.Range("E38").Value = -0.005692860679980649


# Restore sheet protection with the same password.
$ws.Protect("D382")
